$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B-column values (rows 2-127)
$ws.Range("B2").Value = 0.28864
$ws.Range("B3").Value = 0.28864
$ws.Range("B4").Value = 0.27809
$ws.Range("B5").Value = 0.27809
$ws.Range("B6").Value = 0.27809
$ws.Range("B7").Value = 0.27809
$ws.Range("B8").Value = 0.27685
$ws.Range("B9").Value = 0.27809
$ws.Range("B10").Value = 0.27498
$ws.Range("B11").Value = 0.2905
$ws.Range("B12").Value = 0.28305
$ws.Range("B13").Value = 0.28305
$ws.Range("B14").Value = 0.28305
$ws.Range("B15").Value = 0.28305
$ws.Range("B16").Value = 0.27312
$ws.Range("B17").Value = 0.27498
$ws.Range("B18").Value = 0.27685
$ws.Range("B19").Value = 0.26567
$ws.Range("B20").Value = 0.2905
$ws.Range("B21").Value = 0.2725
$ws.Range("B22").Value = 0.28802
$ws.Range("B23").Value = 0.2725
$ws.Range("B24").Value = 0.2843
$ws.Range("B25").Value = 0.28367
$ws.Range("B26").Value = 0.2843
$ws.Range("B27").Value = 0.27312
$ws.Range("B28").Value = 0.27809
$ws.Range("B29").Value = 0.2843
$ws.Range("B30").Value = 0.27809
$ws.Range("B31").Value = 0.26754
$ws.Range("B32").Value = 0.27747
$ws.Range("B33").Value = 0.28181
$ws.Range("B35").Value = 0.28988
$ws.Range("B36").Value = 0.28988
$ws.Range("B37").Value = 0.28988
$ws.Range("B38").Value = 0.29981
$ws.Range("B39").Value = 0.29671
$ws.Range("B40").Value = 0.28988
$ws.Range("B41").Value = 0.28926
$ws.Range("B42").Value = 0.27747
$ws.Range("B43").Value = 0.27312
$ws.Range("B44").Value = 0.27312
$ws.Range("B45").Value = 0.28181
$ws.Range("B46").Value = 0.27623
$ws.Range("B47").Value = 0.27747
$ws.Range("B48").Value = 0.27623
$ws.Range("B49").Value = 0.27623
$ws.Range("B50").Value = 0.26567
$ws.Range("B51").Value = 0.27064
$ws.Range("B52").Value = 0.26878
$ws.Range("B53").Value = 0.26567
$ws.Range("B54").Value = 0.25947
$ws.Range("B55").Value = 0.25885
$ws.Range("B56").Value = 0.26567
$ws.Range("B57").Value = 0.25947
$ws.Range("B58").Value = 0.2545
$ws.Range("B59").Value = 0.25885
$ws.Range("B60").Value = 0.27498
$ws.Range("B61").Value = 0.27126
$ws.Range("B62").Value = 0.26505
$ws.Range("B63").Value = 0.25947
$ws.Range("B64").Value = 0.26505
$ws.Range("B65").Value = 0.26071
$ws.Range("B66").Value = 0.26133
$ws.Range("B67").Value = 0.26133
$ws.Range("B68").Value = 0.27188
$ws.Range("B69").Value = 0.26816
$ws.Range("B70").Value = 0.28367
$ws.Range("B71").Value = 0.2874
$ws.Range("B72").Value = 0.28988
$ws.Range("B73").Value = 0.26691
$ws.Range("B74").Value = 0.28119
$ws.Range("B75").Value = 0.27126
$ws.Range("B76").Value = 0.26878
$ws.Range("B77").Value = 0.28181
$ws.Range("B78").Value = 0.28181
$ws.Range("B79").Value = 0.28367
$ws.Range("B80").Value = 0.30292
$ws.Range("B81").Value = 0.29361
$ws.Range("B82").Value = 0.29671
$ws.Range("B83").Value = 0.30912
$ws.Range("B84").Value = 0.31533
$ws.Range("B85").Value = 0.31844
$ws.Range("B86").Value = 0.30912
$ws.Range("B87").Value = 0.29671
$ws.Range("B88").Value = 0.29671
$ws.Range("B90").Value = 0.2905
$ws.Range("B91").Value = 0.30292
$ws.Range("B92").Value = 0.29671
$ws.Range("B93").Value = 0.29671
$ws.Range("B94").Value = 0.2905
$ws.Range("B95").Value = 0.29361
$ws.Range("B96").Value = 0.28492
$ws.Range("B97").Value = 0.28243
$ws.Range("B98").Value = 0.27188
$ws.Range("B99").Value = 0.26816
$ws.Range("B100").Value = 0.2843
$ws.Range("B101").Value = 0.30602
$ws.Range("B102").Value = 0.30602
$ws.Range("B103").Value = 0.29671
$ws.Range("B104").Value = 0.28802
$ws.Range("B105").Value = 0.2843
$ws.Range("B106").Value = 0.27933
$ws.Range("B107").Value = 0.28926
$ws.Range("B108").Value = 0.2905
$ws.Range("B109").Value = 0.30292
$ws.Range("B110").Value = 0.30912
$ws.Range("B111").Value = 0.29671
$ws.Range("B112").Value = 0.28554
$ws.Range("B113").Value = 0.2905
$ws.Range("B114").Value = 0.31533
$ws.Range("B115").Value = 0.30292
$ws.Range("B116").Value = 0.30912
$ws.Range("B117").Value = 0.29361
$ws.Range("B118").Value = 0.31533
$ws.Range("B119").Value = 0.35878
$ws.Range("B120").Value = 0.39603
$ws.Range("B121").Value = 0.39603
$ws.Range("B122").Value = 0.38361
$ws.Range("B123").Value = 0.38982
$ws.Range("B124").Value = 0.38051
$ws.Range("B125").Value = 0.38051
$ws.Range("B126").Value = 0.37741
$ws.Range("B127").Value = 0.40534

# Add new rows 128-133 (A index values + B values)
$ws.Range("A128").Value = 126
$ws.Range("A129").Value = 127
$ws.Range("A130").Value = 128
$ws.Range("A131").Value = 129
$ws.Range("A132").Value = 130
$ws.Range("A133").Value = 131
$ws.Range("B128").Value = 0.3712
$ws.Range("B129").Value = 0.3712
$ws.Range("B130").Value = 0.41775
$ws.Range("B131").Value = 0.45189
$ws.Range("B132").Value = 0.44569
$ws.Range("B133").Value = 0.4581

# Copy the style (bold font, border, centered/top aligned) from A127 down to the new A cells
$ws.Range("A127").Copy()
$ws.Range("A128:A133").PasteSpecial(-4122)
